$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Insert a new row at the top, pushing existing data down
$ws.Rows.Item(1).Insert()

# Set header values
$ws.Range("A1").Value = "VM"
$ws.Range("B1").Value = "Cloudlet"

# Bold font + thin box border around header cells
$headerRange = $ws.Range("A1:B1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Update selection to F3
$ws.Range("F3").Select()
